# Updated cryptos list on Fri Nov 17 17:13:41 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text (preserving the original
# "no explicit style" appearance), since many of the Price/Volume values
# look numeric (e.g. "0.608", "1.00") and Excel would otherwise coerce
# them into real numbers losing formatting (trailing zeros, thousands
# separators used as decimal group markers, percent sign, padding, etc).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "36.536.03"
Set-TextValue $ws.Range("E2") "  -0.04%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "1.932.26"
Set-TextValue $ws.Range("E3") "  -3.73%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  +0.01%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "243.16"
Set-TextValue $ws.Range("E5") "  -1.71%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.608"
Set-TextValue $ws.Range("E6") "  -3.72%  "

# Row 7 - was USDC, now Solana (swapped with row 8)
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D7") "57.22"
Set-TextValue $ws.Range("E7") "  -8.18%  "

# Row 8 - was Solana, now USDC (swapped with row 7)
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  -0.05%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.366"
Set-TextValue $ws.Range("E9") "  -4.55%  "

# Row 10 - OKB
Set-TextValue $ws.Range("D10") "55.04"
Set-TextValue $ws.Range("E10") "  -3.54%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0827"
Set-TextValue $ws.Range("E11") "  +5.20%  "

# Row 12 - TRON
Set-TextValue $ws.Range("E12") "  -0.59%  "

# Row 13 - was WrappedliquidstakedEther2.0, now Polygon
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Range("D13") "0.820"
Set-TextValue $ws.Range("E13") "  -6.62%  "

# Row 14 - was Polygon, now Avalanche
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D14") "21.44"
Set-TextValue $ws.Range("E14") "  -5.55%  "

# Row 15 - was Avalanche, now WrappedliquidstakedEther2.0
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D15") "2.207.17"
Set-TextValue $ws.Range("E15") "  -4.09%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "13.43"
Set-TextValue $ws.Range("E16") "  -4.54%  "

# Row 17 - Polkadot
Set-TextValue $ws.Range("D17") "5.25"
Set-TextValue $ws.Range("E17") "  -4.91%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "1.908.85"
Set-TextValue $ws.Range("E18") "  -4.97%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "36.492.02"
Set-TextValue $ws.Range("E19") "  -0.05%  "

# Row 20 - Litecoin
Set-TextValue $ws.Range("D20") "69.19"
Set-TextValue $ws.Range("E20") "  -3.75%  "

# Row 21 - ShibaInu
Set-TextValue $ws.Range("D21") "0.0₃0863"
Set-TextValue $ws.Range("E21") "  -1.11%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "228.65"
Set-TextValue $ws.Range("E22") "  -4.12%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "4.99"
Set-TextValue $ws.Range("E23") "  -6.33%  "

# Row 24 - Dai
Set-TextValue $ws.Range("E24") "  -0.09%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.46"
Set-TextValue $ws.Range("E25") "  -2.52%  "

# Row 26 - Toncoin
Set-TextValue $ws.Range("D26") "2.28"
Set-TextValue $ws.Range("E26") "  -2.07%  "

# Row 27 - Cosmos
Set-TextValue $ws.Range("D27") "9.35"
Set-TextValue $ws.Range("E27") "  -6.18%  "

# Row 28 - Monero
Set-TextValue $ws.Range("D28") "162.44"
Set-TextValue $ws.Range("E28") "  +1.89%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "19.29"
Set-TextValue $ws.Range("E29") "  -4.19%  "

# Row 30 - Kaspa
Set-TextValue $ws.Range("D30") "0.123"
Set-TextValue $ws.Range("E30") "  -8.24%  "

# Row 31 - Stellar
Set-TextValue $ws.Range("D31") "0.117"
Set-TextValue $ws.Range("E31") "  -3.19%  "

# Row 32 - ImmutableX
Set-TextValue $ws.Range("D32") "1.14"
Set-TextValue $ws.Range("E32") "  -2.87%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "4.68"
Set-TextValue $ws.Range("E33") "  -6.58%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.0627"
Set-TextValue $ws.Range("E34") "  -0.30%  "

# Row 35 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D35") "4.31"
Set-TextValue $ws.Range("E35") "  -3.99%  "

# Row 36 - BinanceUSD
Set-TextValue $ws.Range("E36") "  -0.11%  "

# Row 37 - was WEMIXToken, now THORChain
$ws.Range("B37").Value = "THORChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D37") "6.04"
Set-TextValue $ws.Range("E37") "  -7.01%  "

# Row 38 - was THORChain, now WEMIXToken
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D38") "1.77"
Set-TextValue $ws.Range("E38") "  -3.10%  "

# Row 39 - LidoDAOToken
Set-TextValue $ws.Range("D39") "2.14"
Set-TextValue $ws.Range("E39") "  -8.51%  "

# Row 40 - RenderToken
Set-TextValue $ws.Range("D40") "2.92"
Set-TextValue $ws.Range("E40") "  -8.50%  "

# Row 41 - Cronos
Set-TextValue $ws.Range("D41") "0.0965"
Set-TextValue $ws.Range("E41") "  -4.10%  "

# Row 42 - was HuobiToken, now TrustWalletToken
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D42") "1.18"
Set-TextValue $ws.Range("E42") "  -6.76%  "

# Row 43 - was TrustWalletToken, now HuobiToken
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D43") "2.85"
Set-TextValue $ws.Range("E43") "  -2.38%  "

# Row 44 - VeChain
Set-TextValue $ws.Range("D44") "0.0208"
Set-TextValue $ws.Range("E44") "  -3.28%  "

# Row 45 - InjectiveProtocol
Set-TextValue $ws.Range("D45") "15.82"
Set-TextValue $ws.Range("E45") "  -5.62%  "

# Row 46 - was Maker, now ARBITRUM
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D46") "1.04"
Set-TextValue $ws.Range("E46") "  -7.51%  "

# Row 47 - was ARBITRUM, now Maker
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D47") "1.339.57"
Set-TextValue $ws.Range("E47") "  -1.74%  "

# Row 48 - Aave
Set-TextValue $ws.Range("D48") "87.42"
Set-TextValue $ws.Range("E48") "  -8.47%  "

# Row 49 - FraxShare
Set-TextValue $ws.Range("D49") "7.25"
Set-TextValue $ws.Range("E49") "  -5.44%  "

# Row 50 - MXToken
Set-TextValue $ws.Range("D50") "2.82"
Set-TextValue $ws.Range("E50") "  -2.20%  "

# Row 51 - MultiversX
Set-TextValue $ws.Range("D51") "45.67"
Set-TextValue $ws.Range("E51") "  +2.41%  "
